$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "39e9dc", 43018.0, "female", 62.0),
    @(3, "664549", 43024.0, "male", 28.0),
    @(4, "b4d8aa", 43025.0, "male", 54.0),
    @(5, "51883d", 43026.0, "male", 57.0),
    @(6, "947e40", 43028.0, "female", 23.0),
    @(7, "9aa197", 43028.0, "female", 66.0),
    @(8, "e4b0a2", 43029.0, "female", 13.0),
    @(9, "af0ac0", 43029.0, "male", 10.0),
    @(10, "185911", 43029.0, "female", 34.0),
    @(11, "601d2e", 43030.0, "male", 11.0),
    @(12, "605322", 43030.0, "female", 23.0),
    @(13, "e399b1", 43031.0, "female", 23.0),
    @(14, "e37897", 43036.0, "female", 9.0),
    @(15, "f658bc", 43036.0, "male", 68.0),
    @(16, "a8e9d8", 43037.0, "female", 37.0),
    @(17, "6f3c15", 43037.0, "male", 13.0),
    @(18, "ae87f7", 43040.0, "male", 25.0),
    @(19, "8c1fc8", 43040.0, "male", 28.0),
    @(20, "8c5776", 43041.0, "female", 7.0),
    @(21, "88526e", 43042.0, "female", 20.0),
    @(22, "1b035f", 43043.0, "male", 29.0),
    @(23, "1efd54", 43043.0, "male", 8.0),
    @(24, "778316", 43043.0, "female", 10.0),
    @(25, "501b8f", 43044.0, "male", 23.0),
    @(26, "525dfa", 43045.0, "female", 10.0),
    @(27, "b5ad13", 43046.0, "female", 21.0),
    @(28, "8bed66", 43047.0, "female", 29.0),
    @(29, "426b6d", 43047.0, "female", 7.0),
    @(30, "a402f6", 43048.0, "male", 26.0),
    @(31, "c2a389", 43049.0, "female", 26.0),
    @(32, "93a3ba", 43049.0, "male", 7.0),
    @(33, "99141e", 43051.0, "male", 15.0),
    @(34, "3b43e6", 43051.0, "male", 28.0),
    @(35, "5eb2b0", 43052.0, "female", 7.0),
    @(36, "3a5305", 43052.0, "male", 29.0),
    @(37, "55a452", 43053.0, "female", 33.0),
    @(38, "b7faf4", 43055.0, "female", 10.0),
    @(39, "4945f2", 43055.0, "male", 19.0),
    @(40, "1875e2", 43055.0, "male", 30.0),
    @(41, "59e66c", 43055.0, "male", 9.0),
    @(42, "933811", 43055.0, "male", 16.0),
    @(43, "bcfad3", 43057.0, "female", 31.0),
    @(44, "99abbe", 43057.0, "male", 14.0),
    @(45, "5ea757", 43058.0, "male", 24.0),
    @(46, "944ba3", 43058.0, "female", 30.0),
    @(47, "95fc1d", 43058.0, "female", 15.0),
    @(48, "5c5c05", 43059.0, "female", 21.0),
    @(49, "d5564e", 43059.0, "female", 43.0),
    @(50, "ac8d9d", 43062.0, "female", 5.0),
    @(51, "fea66b", 43063.0, "male", 18.0),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
